$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update analogs for row 2 (Champion RE10PMC5) to add the Chrysler cross-reference
$ws.Range("F2").Value = "Champion|3032&CHRYSLER|SPRE10PMC5"

# Fix the description for row 2, which previously held boilerplate availability text
$ws.Range("G2").Value = "Свеча зажигания Dodge Interpid, Chrysler Sebring с мотором 2.7 до 2005г. Данная деталь в наличии. Оплата товара за наличный расчет."

# Update price for row 2
$ws.Range("I2").Value = 250

# Update the active selection to match the author's final cursor position
$ws.Range("G5").Select()
